$wb = $excel.ActiveWorkbook

# This script applies cached-value corrections to the Leve profit-tracking
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), matching a scheduled
# market-data refresh: currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H-N) are recomputed per row from freshly pulled prices.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1129
$ws.Range("J17").Value = 1129
$ws.Range("L17").Value = 3387
$ws.Range("N17").Value = -3723
$ws.Range("H86").Value = 3916.9443
$ws.Range("I86").Value = 3909.3635
$ws.Range("K86").Value = 3909.3635
$ws.Range("M86").Value = -2786.3635
$ws.Range("H89").Value = 3916.9443
$ws.Range("I89").Value = 3909.3635
$ws.Range("K89").Value = 19546.8175
$ws.Range("M89").Value = -13930.8175
$ws.Range("H116").Value = 4217.3
$ws.Range("I116").Value = 3699.5
$ws.Range("K116").Value = 3699.5
$ws.Range("M116").Value = -257.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 500997.5
$ws.Range("J6").Value = 1995
$ws.Range("L6").Value = 1995
$ws.Range("N6").Value = -2341
$ws.Range("H45").Value = 911
$ws.Range("I45").Value = 911
$ws.Range("K45").Value = 911
$ws.Range("M45").Value = -534
$ws.Range("H50").Value = 7489.4
$ws.Range("I50").Value = 1224
$ws.Range("J50").Value = 11666.333
$ws.Range("K50").Value = 1224
$ws.Range("L50").Value = 11666.333
$ws.Range("M50").Value = -510
$ws.Range("N50").Value = -13094.333
$ws.Range("H61").Value = 3842.5356
$ws.Range("I61").Value = 3688.12
$ws.Range("J61").Value = 5129.3335
$ws.Range("K61").Value = 3688.12
$ws.Range("L61").Value = 5129.3335
$ws.Range("M61").Value = -3476.12
$ws.Range("N61").Value = -5553.3335
$ws.Range("H110").Value = 948.4
$ws.Range("I110").Value = 914.6667
$ws.Range("K110").Value = 914.6667
$ws.Range("M110").Value = 1130.3333
$ws.Range("H122").Value = 2859.75
$ws.Range("J122").Value = 2765
$ws.Range("L122").Value = 8295
$ws.Range("N122").Value = -13195
$ws.Range("H128").Value = 10000
$ws.Range("J128").Value = 10000
$ws.Range("L128").Value = 10000
$ws.Range("N128").Value = -19960
$ws.Range("H136").Value = 3842.5356
$ws.Range("I136").Value = 3688.12
$ws.Range("J136").Value = 5129.3335
$ws.Range("K136").Value = 11064.36
$ws.Range("L136").Value = 15388.0005
$ws.Range("M136").Value = -8514.360000000001
$ws.Range("N136").Value = -20488.0005

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 973
$ws.Range("I20").Value = 952
$ws.Range("K20").Value = 952
$ws.Range("M20").Value = -705

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1159.5
$ws.Range("I107").Value = 453
$ws.Range("J107").Value = 1748.25
$ws.Range("K107").Value = 453
$ws.Range("L107").Value = 1748.25
$ws.Range("M107").Value = 1467
$ws.Range("N107").Value = -5588.25

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 242.91667
$ws.Range("I11").Value = 201.85715
$ws.Range("K11").Value = 605.5714499999999
$ws.Range("M11").Value = -465.5714499999999
$ws.Range("H87").Value = 5497.25
$ws.Range("I87").Value = 5497.25
$ws.Range("K87").Value = 16491.75
$ws.Range("M87").Value = -15243.75
$ws.Range("H90").Value = 5497.25
$ws.Range("I90").Value = 5497.25
$ws.Range("K90").Value = 49475.25
$ws.Range("M90").Value = -43235.25
$ws.Range("H140").Value = 1666.625
$ws.Range("I140").Value = 1484.7142
$ws.Range("K140").Value = 4454.142599999999
$ws.Range("M140").Value = 725.8574000000008

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1611400
$ws.Range("I21").Value = 1014250
$ws.Range("J21").Value = 4000000
$ws.Range("K21").Value = 1014250
$ws.Range("L21").Value = 4000000
$ws.Range("M21").Value = -1014077
$ws.Range("N21").Value = -4000346
$ws.Range("H29").Value = 2506003.8
$ws.Range("I29").Value = 5005003.5
$ws.Range("J29").Value = 7004
$ws.Range("K29").Value = 5005003.5
$ws.Range("L29").Value = 7004
$ws.Range("M29").Value = -5004713.5
$ws.Range("N29").Value = -7584
$ws.Range("H30").Value = 1611400
$ws.Range("I30").Value = 1014250
$ws.Range("J30").Value = 4000000
$ws.Range("K30").Value = 1014250
$ws.Range("L30").Value = 4000000
$ws.Range("M30").Value = -1014145
$ws.Range("N30").Value = -4000210
$ws.Range("H70").Value = 10027.417
$ws.Range("I70").Value = 6250
$ws.Range("J70").Value = 10499.594
$ws.Range("K70").Value = 6250
$ws.Range("L70").Value = 10499.594
$ws.Range("M70").Value = -5980
$ws.Range("N70").Value = -11039.594
$ws.Range("H73").Value = 10027.417
$ws.Range("I73").Value = 6250
$ws.Range("J73").Value = 10499.594
$ws.Range("K73").Value = 6250
$ws.Range("L73").Value = 10499.594
$ws.Range("M73").Value = -5314
$ws.Range("N73").Value = -12371.594
$ws.Range("H122").Value = 11942.782
$ws.Range("I122").Value = 12284.45
$ws.Range("K122").Value = 36853.35000000001
$ws.Range("M122").Value = -34403.35000000001
$ws.Range("H126").Value = 2804.7273
$ws.Range("I126").Value = 2640.7058
$ws.Range("J126").Value = 3362.4
$ws.Range("K126").Value = 7922.117400000001
$ws.Range("L126").Value = 10087.2
$ws.Range("M126").Value = -5452.117400000001
$ws.Range("N126").Value = -15027.2

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2135.9832
$ws.Range("I22").Value = 1934.0869
$ws.Range("J22").Value = 2264.9722
$ws.Range("K22").Value = 1934.0869
$ws.Range("L22").Value = 2264.9722
$ws.Range("M22").Value = -1639.0869
$ws.Range("N22").Value = -2854.9722
$ws.Range("H27").Value = 2135.9832
$ws.Range("I27").Value = 1934.0869
$ws.Range("J27").Value = 2264.9722
$ws.Range("K27").Value = 1934.0869
$ws.Range("L27").Value = 2264.9722
$ws.Range("M27").Value = -1827.0869
$ws.Range("N27").Value = -2478.9722
$ws.Range("H61").Value = 3985.65
$ws.Range("I61").Value = 1292.1818
$ws.Range("K61").Value = 1292.1818
$ws.Range("M61").Value = -1090.1818
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 3985.65
$ws.Range("I113").Value = 1292.1818
$ws.Range("K113").Value = 1292.1818
$ws.Range("M113").Value = 877.8181999999999
$ws.Range("H122").Value = 3387.75
$ws.Range("J122").Value = 2795.5
$ws.Range("L122").Value = 8386.5
$ws.Range("N122").Value = -13286.5
$ws.Range("H132").Value = 3561.3044
$ws.Range("I132").Value = 3464.1785
$ws.Range("K132").Value = 10392.5355
$ws.Range("M132").Value = -7862.5355

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3584
$ws.Range("J4").Value = 375
$ws.Range("L4").Value = 375
$ws.Range("N4").Value = -601
$ws.Range("H122").Value = 3520.2307
$ws.Range("I122").Value = 928
$ws.Range("K122").Value = 2784
$ws.Range("M122").Value = -334
$ws.Range("H132").Value = 2760.8484
$ws.Range("I132").Value = 2557.682
$ws.Range("J132").Value = 3167.182
$ws.Range("K132").Value = 7673.045999999999
$ws.Range("L132").Value = 9501.545999999998
$ws.Range("M132").Value = -5143.045999999999
$ws.Range("N132").Value = -14561.546
